# Recolor the three "Person" ovals and their three connector arrows in the
# student_network example graph from the dark teal (41817E) to the lighter
# muted teal (94AFAE).
#
# VBA-style RGB(0x94, 0xAF, 0xAE) packed as a single Long (R + G*256 + B*65536).
$newColor = 0x94 + (0xAF * 256) + (0xAE * 65536)

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Ovals filled with the old teal color: "Oval 3", "Oval 5", "Oval 6"
$ovalShapeIndexes = @(1, 3, 4)
foreach ($idx in $ovalShapeIndexes) {
    $shp = $s.Shapes.Item($idx)
    $shp.Fill.ForeColor.RGB = $newColor
}

# Curved connectors whose line uses the old teal color: "Curved Connector 38",
# "Curved Connector 41", "Curved Connector 46"
$connectorShapeIndexes = @(8, 9, 10)
foreach ($idx in $connectorShapeIndexes) {
    $shp = $s.Shapes.Item($idx)
    $shp.Line.ForeColor.RGB = $newColor
}
